$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit permutes the species-occurrence data across rows 2, 3, 4, 7, 8, 9
# (rows 5 and 6 are untouched). Capture the "before" values of every row
# that participates in the permutation first, so that writing the new
# values doesn't clobber data we still need to read.

$cols = @("A","B","D","E","F","G","H","Q","R","Y","AA")
$rows = @(2,3,4,7,8,9)

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        # Value2 (not Value) is needed to reliably read back primitive
        # values through this COM shim.
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Target row <= source row (cyclic permutation of the captured data)
$mapping = @{
    2 = 9
    3 = 8
    4 = 2
    7 = 4
    8 = 3
    9 = 7
}

$dateCols = @("Y","AA")

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $src = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $cell = $ws.Range("$c$targetRow")
        if ($dateCols -contains $c) {
            # Keep these as plain text (e.g. "2023-09-06") instead of
            # letting Excel reinterpret them as date serial numbers.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $src[$c]
    }
}
